$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "eps" (column D) values for rows 2-28
$dValues = @(
    -1.11959252454608,
    -1.09885392870468,
    -1.21451658478844,
    -1.18533373722185,
    -1.06406957393257,
    -1.41067982463927,
    -1.41469986801901,
    -1.32741357648292,
    -1.54527704416407,
    -1.50163549498883,
    -1.40236857634035,
    -1.41255933465945,
    -1.60945924546529,
    -1.60827818722964,
    -1.59673509488722,
    -1.65624235304851,
    -1.46535167676942,
    -1.46918836733134,
    -1.52042277882014,
    -1.58941654194888,
    -1.59601169533397,
    -1.55118418598355,
    -1.71083939446592,
    -1.67952650646475,
    -1.80604592587918,
    -1.69620028585105,
    -1.69373857734667
)

for ($i = 0; $i -lt $dValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $dValues[$i]
}

# Updated "eps_sigma_analytical" (column E) - now a uniform constant for rows 2-28
$ws.Range("E2:E28").Value = 0.124681531161789
